# courseleafPatchControl.xlsx - add cgi step for the search directory
#
# Adds a new "cgi" / "command (chmod)" pair of process rows describing the
# new /web/search/index.cgi script to the CAT and CIM product sheets (and
# to the Include-CGIs template sheet), and renames the "Include-CGIs"
# worksheet to "-Include-CGIs" so the patch engine ignores it (the CGI
# steps previously pulled in from it via the "#include" convention are now
# inlined directly instead).

$wb = $excel.ActiveWorkbook

$wsCAT = $wb.Worksheets.Item("CAT")
$wsCIM = $wb.Worksheets.Item("CIM")
$wsInclude = $wb.Worksheets.Item("Include-CGIs")

# ---------------------------------------------------------------------
# CAT sheet: insert the new rows above the existing "#include" row (17)
# ---------------------------------------------------------------------
$wsCAT.Rows(17).Insert()
$wsCAT.Rows(18).Insert()

$wsCAT.Range("A15:D15").Copy()
$wsCAT.Range("A17:D17").PasteSpecial(-4122)
$wsCAT.Range("A16:D16").Copy()
$wsCAT.Range("A18:D18").PasteSpecial(-4122)
$wsCAT.Range("D16").Copy()
$wsCAT.Range("D17").PasteSpecial(-4122)

$wsCAT.Range("A17").Value = "cgi"
$wsCAT.Range("B17").Value = "ribbit.cgi"
$wsCAT.Range("C17").Value = "/web/search/index.cgi"

$wsCAT.Range("A18").Value = "command"
$wsCAT.Range("B18").Value = "chmod 750 ./web/search/index.cgi"
$wsCAT.Range("C18").Value = "onChangeOnly"

# ---------------------------------------------------------------------
# CIM sheet: the "#include" row (9) is replaced in place by the new cgi
# row, and a new command row is inserted right after it (10)
# ---------------------------------------------------------------------
$wsCIM.Rows(10).Insert()

$wsCIM.Range("A7:D7").Copy()
$wsCIM.Range("A9:D9").PasteSpecial(-4122)
$wsCIM.Range("A8:D8").Copy()
$wsCIM.Range("A10:D10").PasteSpecial(-4122)
$wsCIM.Range("D8").Copy()
$wsCIM.Range("D9").PasteSpecial(-4122)

$wsCIM.Range("A9").Value = "cgi"
$wsCIM.Range("B9").Value = "ribbit.cgi"
$wsCIM.Range("C9").Value = "/web/search/index.cgi"

$wsCIM.Range("A10").Value = "command"
$wsCIM.Range("B10").Value = "chmod 750 ./web/search/index.cgi"
$wsCIM.Range("C10").Value = "onChangeOnly"

# ---------------------------------------------------------------------
# Include-CGIs sheet: insert the new rows above the blank separator row (7)
# ---------------------------------------------------------------------
$wsInclude.Rows(7).Insert()
$wsInclude.Rows(8).Insert()

$wsInclude.Range("A5:D5").Copy()
$wsInclude.Range("A7:D7").PasteSpecial(-4122)
$wsInclude.Range("A6:D6").Copy()
$wsInclude.Range("A8:D8").PasteSpecial(-4122)
$wsInclude.Range("D6").Copy()
$wsInclude.Range("D7").PasteSpecial(-4122)

$wsInclude.Range("A7").Value = "cgi"
$wsInclude.Range("B7").Value = "ribbit.cgi"
$wsInclude.Range("C7").Value = "/web/search/index.cgi"

$wsInclude.Range("A8").Value = "command"
$wsInclude.Range("B8").Value = "chmod 750 ./web/search/index.cgi"
$wsInclude.Range("C8").Value = "onChangeOnly"

# ---------------------------------------------------------------------
# Rename Include-CGIs -> -Include-CGIs so it's now skipped by the patch
# engine (worksheets starting with '-' are ignored)
# ---------------------------------------------------------------------
$wsInclude.Name = "-Include-CGIs"

# ---------------------------------------------------------------------
# View / selection state to match the saved workbook
# ---------------------------------------------------------------------
$wsCAT.Activate()
$wsCAT.Range("A17:XFD18").Select()

$wsInclude.Activate()
$wsInclude.Range("A7:XFD8").Select()

$wsCIM.Activate()
$wsCIM.Range("B18").Select()
